# Updated cryptos list values (Price + Volume(1h)) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Some of the new price strings parse as plain numbers (e.g. "327.89"),
# which would make Excel coerce them into the Number type instead of Text.
# Force those cells to Text first so the stored value keeps its original
# "stringy" representation (matches source data), then drop the temporary
# number-format override with ClearFormats so no stray style is left behind.
$priceTextCells = @("D5", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D19", "D22", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.868.66'
$ws.Range("D3").Value = '1.769.32'
$ws.Range("D5").Value = '327.89'
$ws.Range("D7").Value = '0.4487'
$ws.Range("D8").Value = '0.3555'
$ws.Range("D9").Value = '0.07420'
$ws.Range("D10").Value = '41.98'
$ws.Range("D12").Value = '1.002'
$ws.Range("D13").Value = '20.93'
$ws.Range("D14").Value = '6.024'
$ws.Range("D15").Value = '7.245'
$ws.Range("D16").Value = '1.772.21'
$ws.Range("D17").Value = '93.21'
$ws.Range("D19").Value = '0.06436'
$ws.Range("D22").Value = '5.777'
$ws.Range("D23").Value = '27.901.88'
$ws.Range("D25").Value = '2.107'
$ws.Range("D26").Value = '162.73'
$ws.Range("D27").Value = '20.36'
$ws.Range("D28").Value = '1.977.07'
$ws.Range("D29").Value = '2.157'
$ws.Range("D30").Value = '124.64'
$ws.Range("D31").Value = '1.108'
$ws.Range("D32").Value = '0.09189'
$ws.Range("D33").Value = '5.609'
$ws.Range("D34").Value = '3.658'
$ws.Range("D35").Value = '11.86'
$ws.Range("D36").Value = '0.02290'
$ws.Range("D37").Value = '0.06094'
$ws.Range("D38").Value = '0.2101'
$ws.Range("D39").Value = '0.6327'
$ws.Range("D40").Value = '4.957'
$ws.Range("D41").Value = '1.183'
$ws.Range("D42").Value = '1.399'
$ws.Range("D43").Value = '7.889'
$ws.Range("D44").Value = '13.35'
$ws.Range("D45").Value = '3.741'
$ws.Range("D46").Value = '0.5902'
$ws.Range("D47").Value = '122.36'
$ws.Range("D48").Value = '1.956'
$ws.Range("D49").Value = '0.06905'
$ws.Range("D50").Value = '1.136'
$ws.Range("D51").Value = '72.94'

foreach ($cellRef in $priceTextCells) {
    $ws.Range($cellRef).ClearFormats()
}

# --- Volume(1h) (column E) updates ---
# These are already non-numeric strings (padded with spaces and a % sign)
# so a plain .Value assignment keeps them stored as text.
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("E3").Value = '  +1.85%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("E13").Value = '  +2.92%  '
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("E16").Value = '  +1.83%  '
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("E21").Value = '  +3.05%  '
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").Value = '  +2.47%  '
$ws.Range("E28").Value = '  +2.09%  '
$ws.Range("E29").Value = '  +5.83%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  +6.27%  '
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("E36").Value = '  +1.20%  '
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("E40").Value = '  +1.82%  '
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("E42").Value = '  +2.00%  '
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("E44").Value = '  +2.98%  '
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("E48").Value = '  +1.99%  '
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("E51").Value = '  +2.50%  '
